$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-10-19T23:55:50.192304"
    3 = "2025-10-19T23:55:50.193303"
    4 = "2025-10-19T23:55:50.193303"
    5 = "2025-10-19T23:55:50.193303"
    6 = "2025-10-19T23:55:50.193303"
    7 = "2025-10-19T23:55:50.194306"
    8 = "2025-10-19T23:55:50.194306"
    9 = "2025-10-19T23:55:50.194306"
    10 = "2025-10-19T23:55:50.194306"
    11 = "2025-10-19T23:55:50.195306"
    12 = "2025-10-19T23:55:50.195306"
    13 = "2025-10-19T23:55:50.195306"
    14 = "2025-10-19T23:55:50.195306"
    15 = "2025-10-19T23:55:50.196308"
    16 = "2025-10-19T23:55:50.196308"
    17 = "2025-10-19T23:55:50.196308"
    18 = "2025-10-19T23:55:50.196308"
    19 = "2025-10-19T23:55:50.197308"
    20 = "2025-10-19T23:55:50.197308"
    21 = "2025-10-19T23:55:50.197308"
    22 = "2025-10-19T23:55:50.197308"
    23 = "2025-10-19T23:55:50.197308"
    24 = "2025-10-19T23:55:50.198303"
    25 = "2025-10-19T23:55:50.198303"
    26 = "2025-10-19T23:55:50.198303"
    27 = "2025-10-19T23:55:50.199304"
    28 = "2025-10-19T23:55:50.199304"
    29 = "2025-10-19T23:55:50.199304"
    30 = "2025-10-19T23:55:50.199304"
    31 = "2025-10-19T23:55:50.199304"
    32 = "2025-10-19T23:55:50.202358"
    33 = "2025-10-19T23:55:50.202358"
    34 = "2025-10-19T23:55:50.202894"
    35 = "2025-10-19T23:55:50.202894"
    36 = "2025-10-19T23:55:50.202894"
    37 = "2025-10-19T23:55:50.202894"
    38 = "2025-10-19T23:55:50.202894"
    39 = "2025-10-19T23:55:50.202894"
    40 = "2025-10-19T23:55:50.203884"
    41 = "2025-10-19T23:55:50.203884"
    42 = "2025-10-19T23:55:50.203884"
    43 = "2025-10-19T23:55:50.203884"
    44 = "2025-10-19T23:55:50.203884"
    45 = "2025-10-19T23:55:50.204884"
    46 = "2025-10-19T23:55:50.238141"
    47 = "2025-10-19T23:55:50.238141"
    48 = "2025-10-19T23:55:50.238141"
    49 = "2025-10-19T23:55:50.239142"
    50 = "2025-10-19T23:55:50.239142"
    51 = "2025-10-19T23:55:50.239142"
    52 = "2025-10-19T23:55:50.241163"
    53 = "2025-10-19T23:55:50.242146"
    54 = "2025-10-19T23:55:50.242146"
    55 = "2025-10-19T23:55:50.244144"
    56 = "2025-10-19T23:55:50.245154"
    57 = "2025-10-19T23:55:50.246146"
    58 = "2025-10-19T23:55:50.246146"
    59 = "2025-10-19T23:55:50.246146"
    60 = "2025-10-19T23:55:50.247144"
    61 = "2025-10-19T23:55:50.247144"
    62 = "2025-10-19T23:55:50.247144"
    63 = "2025-10-19T23:55:50.247144"
    64 = "2025-10-19T23:55:50.248144"
    65 = "2025-10-19T23:55:50.248144"
    66 = "2025-10-19T23:55:50.248144"
    67 = "2025-10-19T23:55:50.249143"
    68 = "2025-10-19T23:55:50.249143"
    69 = "2025-10-19T23:55:50.249143"
    70 = "2025-10-19T23:55:50.249143"
    71 = "2025-10-19T23:55:50.250143"
    72 = "2025-10-19T23:55:50.250143"
    73 = "2025-10-19T23:55:50.250143"
    74 = "2025-10-19T23:55:50.251143"
    75 = "2025-10-19T23:55:50.282975"
    76 = "2025-10-19T23:55:50.283974"
    77 = "2025-10-19T23:55:50.283974"
    78 = "2025-10-19T23:55:50.283974"
    79 = "2025-10-19T23:55:50.284993"
    80 = "2025-10-19T23:55:50.284993"
    81 = "2025-10-19T23:55:50.284993"
    82 = "2025-10-19T23:55:50.285977"
    83 = "2025-10-19T23:55:50.285977"
    84 = "2025-10-19T23:55:50.285977"
    85 = "2025-10-19T23:55:50.285977"
    86 = "2025-10-19T23:55:50.286974"
    87 = "2025-10-19T23:55:50.286974"
    88 = "2025-10-19T23:55:50.287975"
    89 = "2025-10-19T23:55:50.287975"
    90 = "2025-10-19T23:55:50.287975"
    91 = "2025-10-19T23:55:50.287975"
    92 = "2025-10-19T23:55:50.288974"
    93 = "2025-10-19T23:55:50.288974"
    94 = "2025-10-19T23:55:50.288974"
    95 = "2025-10-19T23:55:50.288974"
    96 = "2025-10-19T23:55:50.288974"
    97 = "2025-10-19T23:55:50.289971"
    98 = "2025-10-19T23:55:50.289971"
    99 = "2025-10-19T23:55:50.289971"
    100 = "2025-10-19T23:55:50.289971"
    101 = "2025-10-19T23:55:50.289971"
    102 = "2025-10-19T23:55:50.290972"
    103 = "2025-10-19T23:55:50.314190"
    104 = "2025-10-19T23:55:50.314190"
    105 = "2025-10-19T23:55:50.314190"
    106 = "2025-10-19T23:55:50.314190"
    107 = "2025-10-19T23:55:50.315187"
    108 = "2025-10-19T23:55:50.315187"
    109 = "2025-10-19T23:55:50.315187"
    110 = "2025-10-19T23:55:50.315187"
    111 = "2025-10-19T23:55:50.315187"
    112 = "2025-10-19T23:55:50.316183"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
